# Finalizacion - datos de pruebas para los casos CP002/CP003/CP004
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")
$ws.Activate()

# Copy the "plain + border" style (no special alignment) that's already
# used by column A / D2, and paste it onto the cells that are about to
# receive new content so their formatting matches (style index 2).
$ws.Range("D2").Copy()
$ws.Range("D3").PasteSpecial(-4122)
$ws.Range("D4").PasteSpecial(-4122)
$ws.Range("C5").PasteSpecial(-4122)
$ws.Range("D5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new test data (order matters so new shared strings are
# appended in the same sequence as the source workbook).
$ws.Range("D4").Value = "https://www.pcfactory.cl/misdatos"
$ws.Range("C5").Value = "https://www.condorito.com/"
$ws.Range("D3").Value = "RUT o Numero de boleta incorrecto"

# Re-fit columns C and D now that they hold longer text.
$ws.Columns.Item(3).AutoFit()
$ws.Columns.Item(4).AutoFit()

# Leave the same cell selected as in the edited workbook.
$ws.Range("D10").Select()
